# Applies the "Fixing unimportant integration issues" edits:
#  1. Rename sheet "Categories" -> "CodeHierarchies"
#  2. Rename sheet "CategoriesMapping" -> "CodeHierarchiesMapping"
#  3. Update the header/footer font style from "Times New Roman,Normal"
#     to "Times New Roman,Regular" on both of those sheets.
#  4. Update the view state on "CodeHierarchiesMapping" (previously
#     "CategoriesMapping"): scroll position back to A1 (top-left) and the
#     active cell/selection to A26 (previously topLeftCell A22 / A30).

$wb = $excel.ActiveWorkbook

$wsCategories = $wb.Worksheets.Item("Categories")
$wsMapping = $wb.Worksheets.Item("CategoriesMapping")

# --- Rename the sheets ---
$wsCategories.Name = "CodeHierarchies"
$wsMapping.Name = "CodeHierarchiesMapping"

# --- Fix the header/footer font style on both renamed sheets ---
$wsCategories.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$wsCategories.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

$wsMapping.PageSetup.CenterHeader = '&"Times New Roman,Regular"&12&A'
$wsMapping.PageSetup.CenterFooter = '&"Times New Roman,Regular"&12Page &P'

# --- Update the view / selection state on the mapping sheet ---
$wsMapping.Activate()
$wsMapping.Range("A26").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 1
$win.ScrollColumn = 1
